$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 347
$ws.Range("I4").Value = 347
$ws.Range("K4").Value = 347
$ws.Range("M4").Value = -233
$ws.Range("H132").Value = 3416.2654
$ws.Range("I132").Value = 1661.4615
$ws.Range("J132").Value = 10260
$ws.Range("K132").Value = 4984.3845
$ws.Range("L132").Value = 30780
$ws.Range("M132").Value = -2454.3845
$ws.Range("N132").Value = -35840
$ws.Range("H137").Value = 2212.682
$ws.Range("I137").Value = 2204.853
$ws.Range("J137").Value = 2239.3
$ws.Range("K137").Value = 6614.559
$ws.Range("L137").Value = 6717.900000000001
$ws.Range("M137").Value = -4064.559
$ws.Range("N137").Value = -11817.9
$ws.Range("H138").Value = 3649.8635
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 3649.8635
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 10949.5905
$ws.Range("N138").Value = -21229.5905
$ws.Range("M138").Value = $null

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1087.4
$ws.Range("I45").Value = 939.1905
$ws.Range("J45").Value = 1433.2222
$ws.Range("K45").Value = 939.1905
$ws.Range("L45").Value = 1433.2222
$ws.Range("M45").Value = -562.1905
$ws.Range("N45").Value = -2187.2222
$ws.Range("H61").Value = 2978.1555
$ws.Range("I61").Value = 2120.2903
$ws.Range("J61").Value = 4877.7144
$ws.Range("K61").Value = 2120.2903
$ws.Range("L61").Value = 4877.7144
$ws.Range("M61").Value = -1908.2903
$ws.Range("N61").Value = -5301.7144
$ws.Range("H74").Value = 1851.52
$ws.Range("I74").Value = 1159.7949
$ws.Range("J74").Value = 4304
$ws.Range("K74").Value = 1159.7949
$ws.Range("L74").Value = 4304
$ws.Range("M74").Value = -285.7949000000001
$ws.Range("N74").Value = -6052
$ws.Range("H77").Value = 1851.52
$ws.Range("I77").Value = 1159.7949
$ws.Range("J77").Value = 4304
$ws.Range("K77").Value = 5798.9745
$ws.Range("L77").Value = 21520
$ws.Range("M77").Value = -1430.9745
$ws.Range("N77").Value = -30256
$ws.Range("H132").Value = 30301.703
$ws.Range("I132").Value = 43307.04
$ws.Range("J132").Value = 3207.25
$ws.Range("K132").Value = 129921.12
$ws.Range("L132").Value = 9621.75
$ws.Range("M132").Value = -127391.12
$ws.Range("N132").Value = -14681.75
$ws.Range("H136").Value = 2978.1555
$ws.Range("I136").Value = 2120.2903
$ws.Range("J136").Value = 4877.7144
$ws.Range("K136").Value = 6360.8709
$ws.Range("L136").Value = 14633.1432
$ws.Range("M136").Value = -3810.8709
$ws.Range("N136").Value = -19733.1432

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 6296990.5
$ws.Range("I99").Value = 2407927.8
$ws.Range("K99").Value = 2407927.8
$ws.Range("M99").Value = -2406429.8
$ws.Range("H103").Value = 27385.834
$ws.Range("J103").Value = 27385.834
$ws.Range("L103").Value = 27385.834
$ws.Range("N103").Value = -29729.834
$ws.Range("H134").Value = 3329.7632
$ws.Range("I134").Value = 3343.5
$ws.Range("J134").Value = 3300
$ws.Range("K134").Value = 10030.5
$ws.Range("L134").Value = 9900
$ws.Range("M134").Value = -7495.5
$ws.Range("N134").Value = -14970

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2230.7827
$ws.Range("I31").Value = 1286.0605
$ws.Range("J31").Value = 4628.923
$ws.Range("K31").Value = 1286.0605
$ws.Range("L31").Value = 4628.923
$ws.Range("M31").Value = -991.0605
$ws.Range("N31").Value = -5218.923
$ws.Range("H34").Value = 2230.7827
$ws.Range("I34").Value = 1286.0605
$ws.Range("J34").Value = 4628.923
$ws.Range("K34").Value = 1286.0605
$ws.Range("L34").Value = 4628.923
$ws.Range("M34").Value = -1084.0605
$ws.Range("N34").Value = -5032.923
$ws.Range("H58").Value = 1880.0333
$ws.Range("I58").Value = 971.5
$ws.Range("J58").Value = 2675
$ws.Range("K58").Value = 971.5
$ws.Range("L58").Value = 2675
$ws.Range("M58").Value = -768.5
$ws.Range("N58").Value = -3081
$ws.Range("H99").Value = 57803.055
$ws.Range("I99").Value = 144564
$ws.Range("J99").Value = 2591.5454
$ws.Range("K99").Value = 144564
$ws.Range("L99").Value = 2591.5454
$ws.Range("M99").Value = -143066
$ws.Range("N99").Value = -5587.5454
$ws.Range("H122").Value = 2509.1538
$ws.Range("I122").Value = 2592.6365
$ws.Range("J122").Value = 2050
$ws.Range("K122").Value = 7777.9095
$ws.Range("L122").Value = 6150
$ws.Range("M122").Value = -5327.9095
$ws.Range("N122").Value = -11050
$ws.Range("H126").Value = 57803.055
$ws.Range("I126").Value = 144564
$ws.Range("J126").Value = 2591.5454
$ws.Range("K126").Value = 433692
$ws.Range("L126").Value = 7774.6362
$ws.Range("M126").Value = -431222
$ws.Range("N126").Value = -12714.6362
$ws.Range("H132").Value = 1540.9556
$ws.Range("I132").Value = 984.67566
$ws.Range("J132").Value = 4113.75
$ws.Range("K132").Value = 2954.02698
$ws.Range("L132").Value = 12341.25
$ws.Range("M132").Value = -424.0269800000001
$ws.Range("N132").Value = -17401.25
$ws.Range("H134").Value = 1460.5264
$ws.Range("I134").Value = 1050.7333
$ws.Range("J134").Value = 2997.25
$ws.Range("K134").Value = 3152.199900000001
$ws.Range("L134").Value = 8991.75
$ws.Range("M134").Value = -617.1999000000005
$ws.Range("N134").Value = -14061.75
$ws.Range("H136").Value = 1880.0333
$ws.Range("I136").Value = 971.5
$ws.Range("J136").Value = 2675
$ws.Range("K136").Value = 2914.5
$ws.Range("L136").Value = 8025
$ws.Range("M136").Value = -364.5
$ws.Range("N136").Value = -13125

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2656.9143
$ws.Range("I102").Value = 1675
$ws.Range("J102").Value = 4799.273
$ws.Range("K102").Value = 1675
$ws.Range("L102").Value = 4799.273
$ws.Range("M102").Value = -53
$ws.Range("N102").Value = -8043.273
$ws.Range("H132").Value = 3545.0264
$ws.Range("I132").Value = 3346.0688
$ws.Range("J132").Value = 4186.1113
$ws.Range("K132").Value = 10038.2064
$ws.Range("L132").Value = 12558.3339
$ws.Range("M132").Value = -7508.206399999999
$ws.Range("N132").Value = -17618.3339

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 12480.096
$ws.Range("I132").Value = 3831.8333
$ws.Range("J132").Value = 24011.111
$ws.Range("K132").Value = 11495.4999
$ws.Range("L132").Value = 72033.333
$ws.Range("M132").Value = -8965.499899999999
$ws.Range("N132").Value = -77093.333
$ws.Range("H136").Value = 3153.7454
$ws.Range("I136").Value = 1896.9756
$ws.Range("J136").Value = 6834.2856
$ws.Range("K136").Value = 5690.9268
$ws.Range("L136").Value = 20502.8568
$ws.Range("M136").Value = -3140.9268
$ws.Range("N136").Value = -25602.8568

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2008.6571
$ws.Range("I132").Value = 1323
$ws.Range("J132").Value = 3989.4443
$ws.Range("K132").Value = 3969
$ws.Range("L132").Value = 11968.3329
$ws.Range("M132").Value = -1439
$ws.Range("N132").Value = -17028.3329
$ws.Range("H136").Value = 17546118
$ws.Range("I136").Value = 27028690
$ws.Range("J136").Value = 3362.25
$ws.Range("K136").Value = 81086070
$ws.Range("L136").Value = 10086.75
$ws.Range("M136").Value = -81083520
$ws.Range("N136").Value = -15186.75
